# Add a new column D with a header ("Trailing spaces   ", preserving the
# trailing whitespace) plus two numeric data rows, mirroring the existing
# A:C layout. This grows the used range from A1:C3 to A1:D3 and adds a new
# shared string entry for the new header text.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D1").Value = "Trailing spaces   "
$ws.Range("D2").Value = 123
$ws.Range("D3").Value = 456

# Match the selection left behind in the source workbook after the edit.
$ws.Range("E3").Select()
